# Automated hashcode update script.
# Updates the hash values stored in column B of the hashcode sheet,
# matching the rows whose code (column A) is unchanged but whose
# computed hash (column B) changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B89"  = "68439b0181d7876541c13259662e62d3"
    "B99"  = "0dc061740719f94d60c3f2fb1a76b472"
    "B110" = "aa9b18f3904e71ef4dadf88111858b4d"
    "B154" = "d09f757fce10d83c04b40c0872fa2f49"
    "B159" = "567cf77756c9ad1d2efe5d8d378938af"
    "B160" = "e6c7c00b8b6510a2b39353b93e9900e0"
    "B278" = "f628fff06e904e2e47120b72d229abd7"
    "B281" = "01fe8f04b41cfb18b35169e6f343d375"
    "B335" = "11352530e667e1d92b0f0b73ab121312"
    "B339" = "95dbda5d9a8b6ad8dfae2c4599d555fd"
    "B411" = "5016892bf179990dc395a7ae5ebe6a6f"
    "B448" = "c5a9f86f2d2c71529bd5c6e6dee7b713"
    "B542" = "8077e4472b59b3f779e1f9186c1e6d00"
    "B561" = "5860f193406589f9e74d7a8a5210a51e"
    "B574" = "dbf505ff70abf2f98d208c3f09f0e5e2"
    "B592" = "320e9ebd681ed0347b22b3f3e81e84de"
    "B688" = "be0392d0a6d60ca7e9618cc528ca05c4"
    "B693" = "0f57cd6fedeba799532b3b0ba4b4f37e"
    "B711" = "a0988d882b880d8a537e3239a8df8a88"
    "B764" = "b62c6f676d044fe150c53af2473ab411"
    "B776" = "a674c1abc8131bd1104e7863c9f31dd5"
    "B819" = "b6c09b428d120017c20a693b7eebd5f4"
    "B823" = "3ea8f816521a69783b352cf8f251e194"
    "B824" = "6b89c2b53a18291ef70de72b5ff36fa1"
    "B827" = "8984ed957ef45588ab2b7e250414079d"
    "B833" = "3b90ab400a44cba436858271a190263b"
    "B835" = "462b9661f05db7b33cc099b42a4fe747"
    "B863" = "1cc9163b0df7e540f7f31945be6fda33"
    "B877" = "6cb1529c3773a8074f6a97dc67c2a11e"
    "B913" = "e0b86bff7849e4f28182a83327dbfef5"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
